$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2592.238
$ws.Cells.Item(86, 9).Value = 1116.6428
$ws.Cells.Item(86, 10).Value = 5543.4287
$ws.Cells.Item(86, 11).Value = 1116.6428
$ws.Cells.Item(86, 12).Value = 5543.4287
$ws.Cells.Item(86, 13).Value = 6.357199999999921
$ws.Cells.Item(86, 14).Value = -7789.4287
$ws.Cells.Item(89, 8).Value = 2592.238
$ws.Cells.Item(89, 9).Value = 1116.6428
$ws.Cells.Item(89, 10).Value = 5543.4287
$ws.Cells.Item(89, 11).Value = 5583.214
$ws.Cells.Item(89, 12).Value = 27717.1435
$ws.Cells.Item(89, 13).Value = 32.78600000000006
$ws.Cells.Item(89, 14).Value = -38949.14350000001
$ws.Cells.Item(92, 8).Value = 529.64703
$ws.Cells.Item(92, 9).Value = 500.44446
$ws.Cells.Item(92, 10).Value = 562.5
$ws.Cells.Item(92, 11).Value = 500.44446
$ws.Cells.Item(92, 12).Value = 562.5
$ws.Cells.Item(92, 13).Value = 747.5555400000001
$ws.Cells.Item(92, 14).Value = -3058.5
$ws.Cells.Item(100, 8).Value = 15153766
$ws.Cells.Item(100, 9).Value = 23811096
$ws.Cells.Item(100, 11).Value = 23811096
$ws.Cells.Item(100, 13).Value = -23810555
$ws.Cells.Item(137, 8).Value = 1300.1333
$ws.Cells.Item(137, 9).Value = 1040.1
$ws.Cells.Item(137, 10).Value = 1820.2
$ws.Cells.Item(137, 11).Value = 3120.3
$ws.Cells.Item(137, 12).Value = 5460.6
$ws.Cells.Item(137, 13).Value = -570.2999999999997
$ws.Cells.Item(137, 14).Value = -10560.6
$ws.Cells.Item(141, 8).Value = 8200
$ws.Cells.Item(141, 9).Value = 9616.666999999999
$ws.Cells.Item(141, 11).Value = 28850.001
$ws.Cells.Item(141, 13).Value = -23670.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2336.0833
$ws.Cells.Item(2, 9).Value = 2081.5715
$ws.Cells.Item(2, 10).Value = 2692.4
$ws.Cells.Item(2, 11).Value = 2081.5715
$ws.Cells.Item(2, 12).Value = 2692.4
$ws.Cells.Item(2, 13).Value = -1968.5715
$ws.Cells.Item(2, 14).Value = -2918.4
$ws.Cells.Item(61, 8).Value = 942.7857
$ws.Cells.Item(61, 9).Value = 774.7619
$ws.Cells.Item(61, 11).Value = 774.7619
$ws.Cells.Item(61, 13).Value = -562.7619
$ws.Cells.Item(74, 8).Value = 1025.9744
$ws.Cells.Item(74, 9).Value = 870.5417
$ws.Cells.Item(74, 10).Value = 1274.6666
$ws.Cells.Item(74, 11).Value = 870.5417
$ws.Cells.Item(74, 12).Value = 1274.6666
$ws.Cells.Item(74, 13).Value = 3.458300000000008
$ws.Cells.Item(74, 14).Value = -3022.6666
$ws.Cells.Item(77, 8).Value = 1025.9744
$ws.Cells.Item(77, 9).Value = 870.5417
$ws.Cells.Item(77, 10).Value = 1274.6666
$ws.Cells.Item(77, 11).Value = 4352.7085
$ws.Cells.Item(77, 12).Value = 6373.333000000001
$ws.Cells.Item(77, 13).Value = 15.29150000000027
$ws.Cells.Item(77, 14).Value = -15109.333
$ws.Cells.Item(116, 8).Value = 2336.0833
$ws.Cells.Item(116, 9).Value = 2081.5715
$ws.Cells.Item(116, 10).Value = 2692.4
$ws.Cells.Item(116, 11).Value = 2081.5715
$ws.Cells.Item(116, 12).Value = 2692.4
$ws.Cells.Item(116, 13).Value = 212.4285
$ws.Cells.Item(116, 14).Value = -7280.4
$ws.Cells.Item(132, 8).Value = 23280722
$ws.Cells.Item(132, 9).Value = 30304402
$ws.Cells.Item(132, 11).Value = 90913206
$ws.Cells.Item(132, 13).Value = -90910676
$ws.Cells.Item(136, 8).Value = 942.7857
$ws.Cells.Item(136, 9).Value = 774.7619
$ws.Cells.Item(136, 11).Value = 2324.2857
$ws.Cells.Item(136, 13).Value = 225.7143000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2336.0833
$ws.Cells.Item(3, 9).Value = 2081.5715
$ws.Cells.Item(3, 10).Value = 2692.4
$ws.Cells.Item(3, 11).Value = 2081.5715
$ws.Cells.Item(3, 12).Value = 2692.4
$ws.Cells.Item(3, 13).Value = -1967.5715
$ws.Cells.Item(3, 14).Value = -2920.4
$ws.Cells.Item(35, 8).Value = 11250
$ws.Cells.Item(35, 10).Value = 11250
$ws.Cells.Item(35, 12).Value = 11250
$ws.Cells.Item(35, 14).Value = -11870
$ws.Cells.Item(134, 8).Value = 6559.5483
$ws.Cells.Item(134, 9).Value = 2018.2
$ws.Cells.Item(134, 10).Value = 142800
$ws.Cells.Item(134, 11).Value = 6054.6
$ws.Cells.Item(134, 12).Value = 428400
$ws.Cells.Item(134, 13).Value = -3519.6
$ws.Cells.Item(134, 14).Value = -433470

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3092
$ws.Cells.Item(31, 9).Value = 2732.842
$ws.Cells.Item(31, 10).Value = 3660.6667
$ws.Cells.Item(31, 11).Value = 2732.842
$ws.Cells.Item(31, 12).Value = 3660.6667
$ws.Cells.Item(31, 13).Value = -2437.842
$ws.Cells.Item(31, 14).Value = -4250.6667
$ws.Cells.Item(34, 8).Value = 3092
$ws.Cells.Item(34, 9).Value = 2732.842
$ws.Cells.Item(34, 10).Value = 3660.6667
$ws.Cells.Item(34, 11).Value = 2732.842
$ws.Cells.Item(34, 12).Value = 3660.6667
$ws.Cells.Item(34, 13).Value = -2530.842
$ws.Cells.Item(34, 14).Value = -4064.6667
$ws.Cells.Item(58, 8).Value = 1674.6
$ws.Cells.Item(58, 9).Value = 1409.9412
$ws.Cells.Item(58, 10).Value = 2237
$ws.Cells.Item(58, 11).Value = 1409.9412
$ws.Cells.Item(58, 12).Value = 2237
$ws.Cells.Item(58, 13).Value = -1206.9412
$ws.Cells.Item(58, 14).Value = -2643
$ws.Cells.Item(132, 8).Value = 54973.844
$ws.Cells.Item(132, 9).Value = 1576.3077
$ws.Cells.Item(132, 11).Value = 4728.9231
$ws.Cells.Item(132, 13).Value = -2198.9231
$ws.Cells.Item(134, 8).Value = 5639.457
$ws.Cells.Item(134, 9).Value = 905.8276
$ws.Cells.Item(134, 11).Value = 2717.4828
$ws.Cells.Item(134, 13).Value = -182.4827999999998
$ws.Cells.Item(136, 8).Value = 1674.6
$ws.Cells.Item(136, 9).Value = 1409.9412
$ws.Cells.Item(136, 10).Value = 2237
$ws.Cells.Item(136, 11).Value = 4229.8236
$ws.Cells.Item(136, 12).Value = 6711
$ws.Cells.Item(136, 13).Value = -1679.8236
$ws.Cells.Item(136, 14).Value = -11811

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 127414.5
$ws.Cells.Item(132, 9).Value = 1190
$ws.Cells.Item(132, 10).Value = 145446.58
$ws.Cells.Item(132, 11).Value = 3570
$ws.Cells.Item(132, 12).Value = 436339.74
$ws.Cells.Item(132, 13).Value = -1040
$ws.Cells.Item(132, 14).Value = -441399.74

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(41, 8).Value = 1003726.2
$ws.Cells.Item(41, 9).Value = 2501667
$ws.Cells.Item(41, 11).Value = 2501667
$ws.Cells.Item(41, 13).Value = -2501229
$ws.Cells.Item(132, 8).Value = 558774.75
$ws.Cells.Item(132, 9).Value = 1003094.8
$ws.Cells.Item(132, 10).Value = 3374.75
$ws.Cells.Item(132, 11).Value = 3009284.4
$ws.Cells.Item(132, 12).Value = 10124.25
$ws.Cells.Item(132, 13).Value = -3006754.4
$ws.Cells.Item(132, 14).Value = -15184.25
$ws.Cells.Item(136, 8).Value = 8809.549999999999
$ws.Cells.Item(136, 9).Value = 9499.5
$ws.Cells.Item(136, 10).Value = 7774.625
$ws.Cells.Item(136, 11).Value = 28498.5
$ws.Cells.Item(136, 12).Value = 23323.875
$ws.Cells.Item(136, 13).Value = -25948.5
$ws.Cells.Item(136, 14).Value = -28423.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 76967096
$ws.Cells.Item(132, 9).Value = 281252130
$ws.Cells.Item(132, 10).Value = 2681626.5
$ws.Cells.Item(132, 11).Value = 843756390
$ws.Cells.Item(132, 12).Value = 8044879.5
$ws.Cells.Item(132, 13).Value = -843753860
$ws.Cells.Item(132, 14).Value = -8049939.5
$ws.Cells.Item(136, 8).Value = 36727.75
$ws.Cells.Item(136, 9).Value = 46394.184
$ws.Cells.Item(136, 10).Value = 1284.1666
$ws.Cells.Item(136, 11).Value = 139182.552
$ws.Cells.Item(136, 12).Value = 3852.4998
$ws.Cells.Item(136, 13).Value = -136632.552
